$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to remain text, avoiding Excel auto-converting
# numeric-looking strings (e.g. "23.45") into actual numbers.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '27.904.42'
$ws.Range("E2").Value = '  +1.46%  '

Set-TextValue $ws.Range("D3") '1.642.08'
$ws.Range("E3").Value = '  +1.49%  '

Set-TextValue $ws.Range("D5") '212.69'
$ws.Range("E5").Value = '  +0.53%  '

Set-TextValue $ws.Range("D6") '0.525'
$ws.Range("E6").Value = '  +0.62%  '

Set-TextValue $ws.Range("D7") '1.00'
$ws.Range("E7").Value = '  -0.23%  '

Set-TextValue $ws.Range("D8") '23.45'
$ws.Range("E8").Value = '  +1.84%  '

$ws.Range("E9").Value = '  +1.75%  '

$ws.Range("E10").Value = '  +0.33%  '

Set-TextValue $ws.Range("D11") '0.0869'
$ws.Range("E11").Value = '  -1.91%  '

$ws.Range("E12").Value = '  +1.31%  '

Set-TextValue $ws.Range("D13") '1.645.38'
$ws.Range("E13").Value = '  +1.61%  '

$ws.Range("E14").Value = '  +1.07%  '

Set-TextValue $ws.Range("D15") '0.563'
$ws.Range("E15").Value = '  +2.79%  '

$ws.Range("E16").Value = '  +1.99%  '

Set-TextValue $ws.Range("D17") '27.875.46'
$ws.Range("E17").Value = '  +1.34%  '

Set-TextValue $ws.Range("D18") '231.65'
$ws.Range("E18").Value = '  +0.94%  '

Set-TextValue $ws.Range("D19") '7.69'
$ws.Range("E19").Value = '  +1.69%  '

$ws.Range("E20").Value = '  +0.49%  '

$ws.Range("E21").Value = '  -0.25%  '

Set-TextValue $ws.Range("D22") '10.77'
$ws.Range("E22").Value = '  +8.94%  '

Set-TextValue $ws.Range("D23") '4.39'
$ws.Range("E23").Value = '  +2.18%  '

$ws.Range("E24").Value = '  +3.68%  '

Set-TextValue $ws.Range("D25") '151.26'
$ws.Range("E25").Value = '  +1.46%  '

$ws.Range("E26").Value = '  +0.62%  '

$ws.Range("E27").Value = '  +0.53%  '

Set-TextValue $ws.Range("D28") '15.70'
$ws.Range("E28").Value = '  +0.86%  '

$ws.Range("E29").Value = '  -0.23%  '

$ws.Range("E30").Value = '  +0.61%  '

$ws.Range("E31").Value = '  +0.21%  '

$ws.Range("E32").Value = '  +0.91%  '

Set-TextValue $ws.Range("D33") '1.454.16'
$ws.Range("E33").Value = '  +0.17%  '

$ws.Range("E34").Value = '  +1.02%  '

$ws.Range("E35").Value = '  +1.22%  '

Set-TextValue $ws.Range("D37") '0.888'
$ws.Range("E37").Value = '  +2.76%  '

Set-TextValue $ws.Range("D38") '0.566'
$ws.Range("E38").Value = '  +0.75%  '

$ws.Range("E39").Value = '  +0.67%  '

Set-TextValue $ws.Range("D40") '0.917'
$ws.Range("E40").Value = '  -1.45%  '

Set-TextValue $ws.Range("D41") '69.23'
$ws.Range("E41").Value = '  +0.28%  '

$ws.Range("E42").Value = '  -0.17%  '

Set-TextValue $ws.Range("D43") '1.02'
$ws.Range("E43").Value = '  +0.57%  '

Set-TextValue $ws.Range("D44") '2.45'
$ws.Range("E44").Value = '  -1.05%  '

$ws.Range("E45").Value = '  +0.28%  '

$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D46") '5.38'
$ws.Range("E46").Value = '  -0.48%  '

$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D47") '1.78'
$ws.Range("E47").Value = '  +5.93%  '

Set-TextValue $ws.Range("D48") '1.782.89'
$ws.Range("E48").Value = '  +1.05%  '

Set-TextValue $ws.Range("D49") '88.41'
$ws.Range("E49").Value = '  +2.65%  '

$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range("D50") '0.0₆0106'
$ws.Range("E50").Value = '  +0.52%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D51") '0.100'
$ws.Range("E51").Value = '  +1.98%  '
